# Rename worksheets to reflect the "1-1" (v1.1) vocabulary import naming.
#   "Measurement technique" -> "Technique"
#   "Measured property"     -> "Measured property #parameter"
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Measurement technique").Name = "Technique"
$wb.Worksheets.Item("Measured property").Name = "Measured property #parameter"
